$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.611.88"
$ws.Range("E2").Value = "  +3.84%  "
$ws.Range("D3").Value = "2.420.72"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.90"
$ws.Range("E5").Value = "  +4.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.40"
$ws.Range("E6").Value = "  +6.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.512"
$ws.Range("E7").Value = "  +2.39%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  +8.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.45"
$ws.Range("E10").Value = "  +3.71%  "
$ws.Range("E11").Value = "  +1.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.92"
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.94"
$ws.Range("E14").Value = "  +3.10%  "
$ws.Range("D15").Value = "2.799.54"
$ws.Range("E15").Value = "  +2.52%  "
$ws.Range("D16").Value = "2.409.62"
$ws.Range("E16").Value = "  +2.63%  "
$ws.Range("E17").Value = "  +4.51%  "
$ws.Range("D18").Value = "44.465.32"
$ws.Range("E19").Value = "  +3.34%  "
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("E21").Value = "  +3.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.75"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.54"
$ws.Range("E23").Value = "  +2.77%  "
$ws.Range("E24").Value = "  +4.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  +2.42%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.21"
$ws.Range("E27").Value = "  +2.84%  "
$ws.Range("E28").Value = "  -3.71%  "
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.39"
$ws.Range("E30").Value = "  +4.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.50"
$ws.Range("E32").Value = "  +19.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.37"
$ws.Range("E33").Value = "  +10.38%  "
$ws.Range("E34").Value = "  +3.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0776"
$ws.Range("E35").Value = "  +8.43%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("E38").Value = "  +3.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.87"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "120.46"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("E42").Value = "  -2.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.03"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("E44").Value = "  +4.40%  "
$ws.Range("D45").Value = "1.942.21"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("E47").Value = "  +8.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.42"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.68"
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.38"
$ws.Range("E50").Value = "  +6.20%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.02"
$ws.Range("E51").Value = "  +4.90%  "
